# VerveStacks JPN Sets workbook update - 2025-08-05 11:56
# Target sheet: "VEDA_Sets-Proc" (the ~TFM_Psets table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Row 3 (CCGT pset): extend the process-name match pattern with *GasCC*
# and fill in the newly-used SetDesc / AndOr columns.
$ws.Range("B3").Value2 = "ep_gas_combined_cycle*,ep_oil_combined_cycle*,CCGT*,*GasCC*"
$ws.Range("G3").Value2 = "CCGT"
$ws.Range("H3").Value2 = "And"
$ws.Range("I3").Value2 = "Or"

# Row 7 (OCGT/Peaker pset): extend the process-name match pattern with EN*CT*
# and fill in the newly-used AndOr columns.
$ws.Range("B7").Value2 = "ep_gas_gas_turbine*,ep_oil_gas_turbine*,gas turbine*,EN*CT*"
$ws.Range("H7").Value2 = "And"
$ws.Range("I7").Value2 = "Or"

# Row 17 (Nuclear pset): add an exclusion pattern for SMR units plus the
# newly-used AndOr columns.
$ws.Range("B17").Value2 = "-*SMR"
$ws.Range("H17").Value2 = "And"
$ws.Range("I17").Value2 = "Or"
